# Update the "dSF" column (column F) values for the rodón_carlos.xlsx save-data
# sheet to reflect the repulled / recalculated data referenced in the commit
# message ("repull data, push all data, mean calculation").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column F ("dSF")
$updates = @{
    2  = 1
    3  = 4
    5  = 3
    6  = 2
    7  = 5
    9  = -4
    10 = -4
    11 = 3
    12 = 1
    13 = -1
    14 = -2
    15 = 0
    16 = 3
    17 = -1
    19 = -1
    21 = -3
    22 = 3
    23 = 2
    24 = 1
    25 = 8
    26 = 6
    27 = 1
    29 = 2
    30 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
